# Add 12 new data rows (r7:r18) mirroring the existing "Noun" method rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = 42600.782083333332
$ws.Cells.Item(7, 2).Value = "Noun"
$ws.Cells.Item(7, 3).Value = 12311
$ws.Cells.Item(7, 4).Value = 7626
$ws.Cells.Item(7, 5).Value = 1369
$ws.Cells.Item(7, 6).Value = 195
$ws.Cells.Item(7, 7).Value = 74
$ws.Cells.Item(7, 8).Value = 72
$ws.Cells.Item(7, 9).Value = 27
$ws.Cells.Item(7, 10).Value = 3
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 99
$ws.Cells.Item(7, 13).Value = 0

$ws.Cells.Item(8, 1).Value = 42600.804652777777
$ws.Cells.Item(8, 2).Value = "Noun"
$ws.Cells.Item(8, 3).Value = 11465
$ws.Cells.Item(8, 4).Value = 7603
$ws.Cells.Item(8, 5).Value = 1365
$ws.Cells.Item(8, 6).Value = 193
$ws.Cells.Item(8, 7).Value = 70
$ws.Cells.Item(8, 8).Value = 73
$ws.Cells.Item(8, 9).Value = 26
$ws.Cells.Item(8, 10).Value = 3
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 99
$ws.Cells.Item(8, 13).Value = 0

$ws.Cells.Item(9, 1).Value = 42600.806793981479
$ws.Cells.Item(9, 2).Value = "Noun"
$ws.Cells.Item(9, 3).Value = 10729
$ws.Cells.Item(9, 4).Value = 7609
$ws.Cells.Item(9, 5).Value = 1367
$ws.Cells.Item(9, 6).Value = 165
$ws.Cells.Item(9, 7).Value = 70
$ws.Cells.Item(9, 8).Value = 69
$ws.Cells.Item(9, 9).Value = 29
$ws.Cells.Item(9, 10).Value = 3
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 99
$ws.Cells.Item(9, 13).Value = 0

$ws.Cells.Item(10, 1).Value = 42600.823240740741
$ws.Cells.Item(10, 2).Value = "Noun"
$ws.Cells.Item(10, 3).Value = 11324
$ws.Cells.Item(10, 4).Value = 7260
$ws.Cells.Item(10, 5).Value = 1315
$ws.Cells.Item(10, 6).Value = 153
$ws.Cells.Item(10, 7).Value = 67
$ws.Cells.Item(10, 8).Value = 68
$ws.Cells.Item(10, 9).Value = 30
$ws.Cells.Item(10, 10).Value = 3
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 99
$ws.Cells.Item(10, 13).Value = 0

$ws.Cells.Item(11, 1).Value = 42600.830231481479
$ws.Cells.Item(11, 2).Value = "Noun"
$ws.Cells.Item(11, 3).Value = 11895
$ws.Cells.Item(11, 4).Value = 7569
$ws.Cells.Item(11, 5).Value = 1366
$ws.Cells.Item(11, 6).Value = 160
$ws.Cells.Item(11, 7).Value = 67
$ws.Cells.Item(11, 8).Value = 70
$ws.Cells.Item(11, 9).Value = 29
$ws.Cells.Item(11, 10).Value = 3
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 99
$ws.Cells.Item(11, 13).Value = 0

$ws.Cells.Item(12, 1).Value = 42600.841689814813
$ws.Cells.Item(12, 2).Value = "Noun"
$ws.Cells.Item(12, 3).Value = 10605
$ws.Cells.Item(12, 4).Value = 7461
$ws.Cells.Item(12, 5).Value = 1302
$ws.Cells.Item(12, 6).Value = 161
$ws.Cells.Item(12, 7).Value = 71
$ws.Cells.Item(12, 8).Value = 69
$ws.Cells.Item(12, 9).Value = 30
$ws.Cells.Item(12, 10).Value = 3
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 99
$ws.Cells.Item(12, 13).Value = 0

$ws.Cells.Item(13, 1).Value = 42600.861238425925
$ws.Cells.Item(13, 2).Value = "Noun"
$ws.Cells.Item(13, 3).Value = 11105
$ws.Cells.Item(13, 4).Value = 6885
$ws.Cells.Item(13, 5).Value = 1201
$ws.Cells.Item(13, 6).Value = 154
$ws.Cells.Item(13, 7).Value = 62
$ws.Cells.Item(13, 8).Value = 70
$ws.Cells.Item(13, 9).Value = 28
$ws.Cells.Item(13, 10).Value = 3
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 99
$ws.Cells.Item(13, 13).Value = 0

$ws.Cells.Item(14, 1).Value = 42600.868344907409
$ws.Cells.Item(14, 2).Value = "Noun"
$ws.Cells.Item(14, 3).Value = 10425
$ws.Cells.Item(14, 4).Value = 6671
$ws.Cells.Item(14, 5).Value = 1156
$ws.Cells.Item(14, 6).Value = 149
$ws.Cells.Item(14, 7).Value = 62
$ws.Cells.Item(14, 8).Value = 70
$ws.Cells.Item(14, 9).Value = 29
$ws.Cells.Item(14, 10).Value = 3
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 99
$ws.Cells.Item(14, 13).Value = 0

$ws.Cells.Item(15, 1).Value = 42600.878541666665
$ws.Cells.Item(15, 2).Value = "Noun"
$ws.Cells.Item(15, 3).Value = 10149
$ws.Cells.Item(15, 4).Value = 6431
$ws.Cells.Item(15, 5).Value = 1079
$ws.Cells.Item(15, 6).Value = 149
$ws.Cells.Item(15, 7).Value = 61
$ws.Cells.Item(15, 8).Value = 70
$ws.Cells.Item(15, 9).Value = 28
$ws.Cells.Item(15, 10).Value = 3
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 99
$ws.Cells.Item(15, 13).Value = 0

$ws.Cells.Item(16, 1).Value = 42600.883969907409
$ws.Cells.Item(16, 2).Value = "Noun"
$ws.Cells.Item(16, 3).Value = 10437
$ws.Cells.Item(16, 4).Value = 6387
$ws.Cells.Item(16, 5).Value = 1094
$ws.Cells.Item(16, 6).Value = 148
$ws.Cells.Item(16, 7).Value = 59
$ws.Cells.Item(16, 8).Value = 71
$ws.Cells.Item(16, 9).Value = 28
$ws.Cells.Item(16, 10).Value = 3
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 99
$ws.Cells.Item(16, 13).Value = 0

$ws.Cells.Item(17, 1).Value = 42600.88517361111
$ws.Cells.Item(17, 2).Value = "Noun"
$ws.Cells.Item(17, 3).Value = 10643
$ws.Cells.Item(17, 4).Value = 6397
$ws.Cells.Item(17, 5).Value = 1093
$ws.Cells.Item(17, 6).Value = 148
$ws.Cells.Item(17, 7).Value = 59
$ws.Cells.Item(17, 8).Value = 71
$ws.Cells.Item(17, 9).Value = 28
$ws.Cells.Item(17, 10).Value = 3
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 99
$ws.Cells.Item(17, 13).Value = 0

$ws.Cells.Item(18, 1).Value = 42600.886354166665
$ws.Cells.Item(18, 2).Value = "Noun"
$ws.Cells.Item(18, 3).Value = 10442
$ws.Cells.Item(18, 4).Value = 6114
$ws.Cells.Item(18, 5).Value = 1039
$ws.Cells.Item(18, 6).Value = 141
$ws.Cells.Item(18, 7).Value = 51
$ws.Cells.Item(18, 8).Value = 73
$ws.Cells.Item(18, 9).Value = 26
$ws.Cells.Item(18, 10).Value = 3
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 99
$ws.Cells.Item(18, 13).Value = 0

# Column A widened slightly (Excel re-ran best-fit after the new, longer date values were added)
$ws.Columns.Item(1).ColumnWidth = 14
